$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# H5: new text "Resumen de Reunión 18" (new shared string)
$ws.Range("H5").Value = "Resumen de Reunión 18"

# Row height adjustments
$ws.Rows.Item(14).RowHeight = 30
$ws.Rows.Item(15).RowHeight = 45
$ws.Rows.Item(16).RowHeight = 30
$ws.Rows.Item(17).RowHeight = 45

# F17: new text "Informe Final de SQA" with wrap text style
$ws.Range("F17").Value = "Informe Final de SQA"
$ws.Range("F17").WrapText = $true

# Selection change to H6
[void]$ws.Range("H6").Select()
